$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header value for the new 2022 column (P)
$ws.Range("P4").Value = 2022

# Data values for column P, rows 5-14 (mirroring column O)
$ws.Range("P5").Value = 96.969944810665083
$ws.Range("P6").Value = 96.173557859042035
$ws.Range("P7").Value = 62.289845326160055
$ws.Range("P8").Value = 100
$ws.Range("P9").Value = 100
$ws.Range("P10").Value = "-"
$ws.Range("P11").Value = 100
$ws.Range("P12").Value = 58.090784503861151
$ws.Range("P13").Value = 100
$ws.Range("P14").Value = 100

# Copy styles from column O to column P so formatting matches (cell by cell
# to avoid the destination range picking up a single blended style)
foreach ($r in 3..14) {
    $ws.Range("O$r").Copy()
    $ws.Range("P$r").PasteSpecial(-4122)  # xlPasteFormats
}

# Update selection to match the target diff
$ws.Range("Q4").Select()
